$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare new row 21 formatting (copy from row 20, which has the same
#     per-column styles we need: s=1 on A/B/C, s=9 on D) ---
$ws.Range("A20:D20").Copy()
$ws.Range("A21:D21").PasteSpecial(-4122)

# --- Update status column (D9:D21) IN PROGRESS -> COMPLETED ---
# D9 is the first cell touched overall, so "COMPLETED" becomes the first
# brand-new shared string appended.
$ws.Range("D9").Value = "COMPLETED"
$ws.Range("D10").Value = "COMPLETED"
$ws.Range("D11").Value = "COMPLETED"
$ws.Range("D12").Value = "COMPLETED"
$ws.Range("D13").Value = "COMPLETED"
$ws.Range("D14").Value = "COMPLETED"
$ws.Range("D15").Value = "COMPLETED"
$ws.Range("D16").Value = "COMPLETED"
$ws.Range("D17").Value = "COMPLETED"
$ws.Range("D18").Value = "COMPLETED"
$ws.Range("D19").Value = "COMPLETED"
$ws.Range("D20").Value = "COMPLETED"
$ws.Range("D21").Value = "COMPLETED"

# --- Row 11: Action Item text changes (second brand-new shared string) ---
$ws.Range("A11").Value = "Player Class (function that determines whether the guessed number is closer to the generated number)"

# --- Row 21 (new row): Action Item / Item ID / Team Member ---
# "UI and the overall design " must be introduced before the row 17/18/19
# replacement text below, to match the target shared-string order.
$ws.Range("A21").Value = "UI and the overall design "
$ws.Range("B21").Value = 13
$ws.Range("C21").Value = "Ming"

# --- Row 18: Action Item text change ---
$ws.Range("A18").Value = "Guessed Number Tracker (keeps track of the number of times a player has guessed)"

# --- Row 17: Action Item text change ---
$ws.Range("A17").Value = "Handle what to do when players quit"

# --- Row 19: Action Item text change ---
$ws.Range("A19").Value = "Win conditions "

# --- Team member reassignments (Column C), rows 9-21 ---
$ws.Range("C9").Value = "Ming"
$ws.Range("C10").Value = "Vishal"
$ws.Range("C11").Value = "Vishal"
$ws.Range("C12").Value = "Vishal"
$ws.Range("C13").Value = "Ashley"
$ws.Range("C14").Value = "Eldin"
$ws.Range("C15").Value = "Eldin"
$ws.Range("C16").Value = "Eldin"
$ws.Range("C17").Value = "Ashley"
$ws.Range("C18").Value = "Eldin"
$ws.Range("C19").Value = "Ashley"
$ws.Range("C20").Value = "Ming"

# --- Server/Client column (G), rows 10-12 change from "Server" label cell
#     content (unchanged text, already "Server"/"Client" - kept for clarity) ---
$ws.Range("G10").Value = "Server"
$ws.Range("G11").Value = "Server"
$ws.Range("G12").Value = "Server"

# --- Update sheet selection to match the edited workbook (G10) ---
$ws.Range("G10").Select()
